# Refresh currentAveragePrice / Leve-profit market-data columns (H:N) on the
# per-class Leve sheets, per the scheduled market-data runner.
# Column layout: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ,
# N=LeveProfitHQ.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: Enchanted Silver Ink
$ws.Range("H28").Value = 1492.75
$ws.Range("I28").Value = 323.66666
$ws.Range("J28").Value = 5000
$ws.Range("K28").Value = 323.66666
$ws.Range("L28").Value = 5000
$ws.Range("M28").Value = 161.33334
$ws.Range("N28").Value = -5970

# Row 29: Weak Blinding Potion
$ws.Range("H29").Value = 1500
$ws.Range("I29").Value = 333.33334
$ws.Range("K29").Value = 1000.00002
$ws.Range("M29").Value = -719.0000200000001

# Row 38: Hi-Potion of Strength
$ws.Range("H38").Value = 7482.5
$ws.Range("I38").Value = 7448
$ws.Range("K38").Value = 22344
$ws.Range("M38").Value = -21972

# Row 40: Horn Glue
$ws.Range("H40").Value = 2737.6
$ws.Range("J40").Value = 2922
$ws.Range("L40").Value = 2922
$ws.Range("N40").Value = -3272

# Row 112: Superior Spiritbond Potion
$ws.Range("H112").Value = 1500
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1500
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4500
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -6716

# Row 137: Magnesia Whetstone
$ws.Range("H137").Value = 3489.8572
$ws.Range("I137").Value = 3404.8333
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 10214.4999
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -7664.499899999999
$ws.Range("N137").Value = -17100

$ws = $wb.Worksheets.Item("ARM")
# Row 74: Titanium Nugget
$ws.Range("H74").Value = 1474
$ws.Range("I74").Value = 1474
$ws.Range("K74").Value = 1474
$ws.Range("M74").Value = -600

# Row 77: Titanium Nugget
$ws.Range("H77").Value = 1474
$ws.Range("I77").Value = 1474
$ws.Range("K77").Value = 7370
$ws.Range("M77").Value = -3002

# Row 102: Tama-hagane Ingot
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 622
$ws.Range("N102").ClearContents()

# Row 104: Molybdenum Kite Shield
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Row 122: High Durium Nugget
$ws.Range("H122").Value = 5647
$ws.Range("I122").Value = 7470.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 22411.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -19961.5
$ws.Range("N122").Value = -10900

# Row 126: Bismuth Ingot
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# Row 132: Mountain Chromite Ingot
$ws.Range("H132").Value = 1200
$ws.Range("I132").Value = 1200
$ws.Range("K132").Value = 3600
$ws.Range("M132").Value = -1070

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Oroshigane Ingot
$ws.Range("H99").Value = 1500
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 1500
$ws.Range("M99").Value = -2

# Row 105: Molybdenum Ingot
$ws.Range("H105").Value = 9924.75
$ws.Range("I105").Value = 9933
$ws.Range("K105").Value = 9933
$ws.Range("M105").Value = -8186

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Walnut Lumber
$ws.Range("H31").Value = 2030.75
$ws.Range("I31").Value = 1874.3334
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 1874.3334
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = -1579.3334
$ws.Range("N31").Value = -3090

# Row 34: Walnut Lumber
$ws.Range("H34").Value = 2030.75
$ws.Range("I34").Value = 1874.3334
$ws.Range("J34").Value = 2500
$ws.Range("K34").Value = 1874.3334
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = -1672.3334
$ws.Range("N34").Value = -2904

# Row 104: Zelkova Necklace
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Row 105: Zelkova Lumber
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").ClearContents()

# Row 132: Ginseng Lumber
$ws.Range("H132").Value = 8444.444
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

# Row 134: Ceiba Lumber
$ws.Range("H134").Value = 2910.1875
$ws.Range("I134").Value = 1570.5
$ws.Range("K134").Value = 4711.5
$ws.Range("M134").Value = -2176.5

$ws = $wb.Worksheets.Item("CUL")
# Row 4: Boiled Egg
$ws.Range("H4").Value = 497.33334
$ws.Range("J4").Value = 494
$ws.Range("L4").Value = 1482
$ws.Range("N4").Value = -1706

# Row 6: Marmot Steak
$ws.Range("H6").Value = 45.25
$ws.Range("J6").Value = 100
$ws.Range("L6").Value = 300
$ws.Range("N6").Value = -526

# Row 55: Pastry Fish
$ws.Range("H55").Value = 2719.8
$ws.Range("J55").Value = 3199.75
$ws.Range("L55").Value = 9599.25
$ws.Range("N55").Value = -9953.25

# Row 98: Rice Vinegar
$ws.Range("H98").Value = 449.83334
$ws.Range("J98").Value = 349.66666
$ws.Range("L98").Value = 1048.99998
$ws.Range("N98").Value = -4044.99998

$ws = $wb.Worksheets.Item("GSM")
# Row 95: Koppranickel Temple Chain
$ws.Range("H95").Value = 35391
$ws.Range("J95").Value = 35391
$ws.Range("L95").Value = 35391
$ws.Range("N95").Value = -40883

# Row 113: Manasilver Nugget
$ws.Range("H113").Value = 3472.8
$ws.Range("I113").Value = 3719.7778
$ws.Range("J113").Value = 1250
$ws.Range("K113").Value = 3719.7778
$ws.Range("L113").Value = 1250
$ws.Range("M113").Value = -1549.7778
$ws.Range("N113").Value = -5590

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Aldgoat Leather
$ws.Range("H22").Value = 2900.5
$ws.Range("I22").Value = 2900.5
$ws.Range("K22").Value = 2900.5
$ws.Range("M22").Value = -2605.5

# Row 27: Aldgoat Leather
$ws.Range("H27").Value = 2900.5
$ws.Range("I27").Value = 2900.5
$ws.Range("K27").Value = 2900.5
$ws.Range("M27").Value = -2793.5

# Row 55: Peiste Leather
$ws.Range("H55").Value = 1975.4
$ws.Range("I55").Value = 1958.3334
$ws.Range("J55").Value = 2001
$ws.Range("K55").Value = 1958.3334
$ws.Range("L55").Value = 2001
$ws.Range("M55").Value = -1785.3334
$ws.Range("N55").Value = -2347

# Row 61: Raptor Leather
$ws.Range("H61").Value = 2664.6667
$ws.Range("I61").Value = 3499.5
$ws.Range("J61").Value = 995
$ws.Range("K61").Value = 3499.5
$ws.Range("L61").Value = 995
$ws.Range("M61").Value = -3297.5
$ws.Range("N61").Value = -1399

# Row 68: Wyvern Leather
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

# Row 71: Wyvern Leather
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

# Row 93: Gagana Leather
$ws.Range("H93").Value = 4949.5
$ws.Range("I93").Value = 4966
$ws.Range("K93").Value = 4966
$ws.Range("M93").Value = -3718

# Row 113: Atrociraptor Leather
$ws.Range("H113").Value = 2664.6667
$ws.Range("I113").Value = 3499.5
$ws.Range("J113").Value = 995
$ws.Range("K113").Value = 3499.5
$ws.Range("L113").Value = 995
$ws.Range("M113").Value = -1329.5
$ws.Range("N113").Value = -5335

# Row 122: Gaja Leather
$ws.Range("H122").Value = 6500
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -34900

# Row 132: Silver Lobo Leather
$ws.Range("H132").Value = 11726.182
$ws.Range("I132").Value = 11898.8
$ws.Range("K132").Value = 35696.39999999999
$ws.Range("M132").Value = -33166.39999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 136: Sarcenet Cloth
$ws.Range("H136").Value = 1980.75
$ws.Range("I136").Value = 1549.4286
$ws.Range("K136").Value = 4648.2858
$ws.Range("M136").Value = -2098.2858
